$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 545.7742
$arr[0,1] = 383.46667
$arr[0,2] = 697.9375
$arr[0,3] = 1150.40001
$arr[0,4] = 2093.8125
$arr[0,5] = -152.4000100000001
$arr[0,6] = -4089.8125
$ws.Range("H80:N80").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 545.7742
$arr[0,1] = 383.46667
$arr[0,2] = 697.9375
$arr[0,3] = 3451.20003
$arr[0,4] = 6281.4375
$arr[0,5] = 1540.79997
$arr[0,6] = -16265.4375
$ws.Range("H83:N83").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1320.6471
$arr[0,1] = 665.3333
$arr[0,2] = 6235.5
$arr[0,3] = 665.3333
$arr[0,4] = 6235.5
$arr[0,5] = 832.6667
$arr[0,6] = -9231.5
$ws.Range("H98:N98").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3093.9443
$arr[0,1] = 1999
$arr[0,2] = 3114.6038
$arr[0,3] = 5997
$arr[0,4] = 9343.811399999999
$arr[0,5] = -4889
$arr[0,6] = -11559.8114
$ws.Range("H112:N112").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1320.6471
$arr[0,1] = 665.3333
$arr[0,2] = 6235.5
$arr[0,3] = 1995.9999
$arr[0,4] = 18706.5
$arr[0,5] = 454.0001
$arr[0,6] = -23606.5
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5091.409
$arr[0,1] = 4632.1577
$arr[0,2] = 8000
$arr[0,3] = 13896.4731
$arr[0,4] = 24000
$arr[0,5] = -11366.4731
$arr[0,6] = -29060
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1715.6
$arr[0,1] = 1508.5518
$arr[0,2] = 2716.3333
$arr[0,3] = 1508.5518
$arr[0,4] = 2716.3333
$arr[0,5] = -1395.5518
$arr[0,6] = -2942.3333
$ws.Range("H2:N2").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 9234.799999999999
$arr[0,1] = 4918.5415
$arr[0,2] = 26499.834
$arr[0,3] = 4918.5415
$arr[0,4] = 26499.834
$arr[0,5] = -4631.5415
$arr[0,6] = -27073.834
$ws.Range("H32:N32").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5555.6924
$arr[0,1] = 3699.25
$arr[0,2] = 11743.833
$arr[0,3] = 3699.25
$arr[0,4] = 11743.833
$arr[0,5] = -3487.25
$arr[0,6] = -12167.833
$ws.Range("H61:N61").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8070.5
$arr[0,1] = 3247
$arr[0,2] = 9999.9
$arr[0,3] = 3247
$arr[0,4] = 9999.9
$arr[0,5] = -2561
$arr[0,6] = -11371.9
$ws.Range("H63:N63").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8070.5
$arr[0,1] = 3247
$arr[0,2] = 9999.9
$arr[0,3] = 16235
$arr[0,4] = 49999.5
$arr[0,5] = -12803
$arr[0,6] = -56863.5
$ws.Range("H66:N66").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2344.65
$arr[0,1] = 2470.9375
$arr[0,2] = 1839.5
$arr[0,3] = 2470.9375
$arr[0,4] = 1839.5
$arr[0,5] = -1596.9375
$arr[0,6] = -3587.5
$ws.Range("H74:N74").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2344.65
$arr[0,1] = 2470.9375
$arr[0,2] = 1839.5
$arr[0,3] = 12354.6875
$arr[0,4] = 9197.5
$arr[0,5] = -7986.6875
$arr[0,6] = -17933.5
$ws.Range("H77:N77").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 814.1667
$arr[0,1] = 824.5454999999999
$arr[0,2] = 700
$arr[0,3] = 824.5454999999999
$arr[0,4] = 700
$arr[0,5] = -328.5454999999999
$arr[0,6] = -1692
$ws.Range("H97:N97").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 12328.277
$arr[0,1] = 13032.883
$arr[0,2] = 350
$arr[0,3] = 13032.883
$arr[0,4] = 350
$arr[0,5] = -11410.883
$arr[0,6] = -3594
$ws.Range("H102:N102").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 31746.75
$arr[0,1] = 0
$arr[0,2] = 31746.75
$arr[0,3] = 0
$arr[0,4] = 31746.75
$arr[0,5] = ""
$arr[0,6] = -34700.75
$ws.Range("H112:N112").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1715.6
$arr[0,1] = 1508.5518
$arr[0,2] = 2716.3333
$arr[0,3] = 1508.5518
$arr[0,4] = 2716.3333
$arr[0,5] = 785.4482
$arr[0,6] = -7304.3333
$ws.Range("H116:N116").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5555.6924
$arr[0,1] = 3699.25
$arr[0,2] = 11743.833
$arr[0,3] = 11097.75
$arr[0,4] = 35231.499
$arr[0,5] = -8547.75
$arr[0,6] = -40331.499
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1715.6
$arr[0,1] = 1508.5518
$arr[0,2] = 2716.3333
$arr[0,3] = 1508.5518
$arr[0,4] = 2716.3333
$arr[0,5] = -1394.5518
$arr[0,6] = -2944.3333
$ws.Range("H3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5548
$arr[0,1] = 5175.6665
$arr[0,2] = 5994.8
$arr[0,3] = 5175.6665
$arr[0,4] = 5994.8
$arr[0,5] = -4928.6665
$arr[0,6] = -6488.8
$ws.Range("H20:N20").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4839.4
$arr[0,1] = 4839.4
$arr[0,2] = 0
$arr[0,3] = 4839.4
$arr[0,4] = 0
$arr[0,5] = -3341.4
$arr[0,6] = ""
$ws.Range("H99:N99").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2810.5715
$arr[0,1] = 2316.6667
$arr[0,2] = 5774
$arr[0,3] = 2316.6667
$arr[0,4] = 5774
$arr[0,5] = -569.6667000000002
$arr[0,6] = -9268
$ws.Range("H105:N105").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3956.5715
$arr[0,1] = 7100
$arr[0,2] = 2699.2
$arr[0,3] = 7100
$arr[0,4] = 2699.2
$arr[0,5] = -6897
$arr[0,6] = -3105.2
$ws.Range("H58:N58").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6553.067
$arr[0,1] = 6867.3335
$arr[0,2] = 6081.6665
$arr[0,3] = 6867.3335
$arr[0,4] = 6081.6665
$arr[0,5] = -6243.3335
$arr[0,6] = -7329.6665
$ws.Range("H62:N62").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6553.067
$arr[0,1] = 6867.3335
$arr[0,2] = 6081.6665
$arr[0,3] = 34336.6675
$arr[0,4] = 30408.3325
$arr[0,5] = -31216.6675
$arr[0,6] = -36648.3325
$ws.Range("H65:N65").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 27799942
$arr[0,1] = 27799942
$arr[0,2] = 0
$arr[0,3] = 27799942
$arr[0,4] = 0
$arr[0,5] = -27798819
$arr[0,6] = ""
$ws.Range("H86:N86").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 27799942
$arr[0,1] = 27799942
$arr[0,2] = 0
$arr[0,3] = 138999710
$arr[0,4] = 0
$arr[0,5] = -138994094
$arr[0,6] = ""
$ws.Range("H89:N89").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 70000
$arr[0,1] = 0
$arr[0,2] = 70000
$arr[0,3] = 0
$arr[0,4] = 70000
$arr[0,5] = ""
$arr[0,6] = -75242
$ws.Range("H104:N104").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3631.087
$arr[0,1] = 3896.7646
$arr[0,2] = 2878.3333
$arr[0,3] = 11690.2938
$arr[0,4] = 8634.999899999999
$arr[0,5] = -9240.293799999999
$arr[0,6] = -13534.9999
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3956.5715
$arr[0,1] = 7100
$arr[0,2] = 2699.2
$arr[0,3] = 21300
$arr[0,4] = 8097.599999999999
$arr[0,5] = -18750
$arr[0,6] = -13197.6
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 212338080
$arr[0,1] = 333563460
$arr[0,2] = 30500000
$arr[0,3] = 1000690380
$arr[0,4] = 91500000
$arr[0,5] = -1000690268
$arr[0,6] = -91500224
$ws.Range("H4:N4").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7999.857
$arr[0,1] = 7999
$arr[0,2] = 8000
$arr[0,3] = 23997
$arr[0,4] = 24000
$arr[0,5] = -23841
$arr[0,6] = -24312
$ws.Range("H49:N49").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5099.25
$arr[0,1] = 3799
$arr[0,2] = 9000
$arr[0,3] = 11397
$arr[0,4] = 27000
$arr[0,5] = -10154
$arr[0,6] = -29486
$ws.Range("H118:N118").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = ""
$arr[0,6] = ""
$ws.Range("H104:N104").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 74999
$arr[0,1] = 0
$arr[0,2] = 74999
$arr[0,3] = 0
$arr[0,4] = 74999
$arr[0,5] = ""
$arr[0,6] = -85359
$ws.Range("H141:N141").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10999
$arr[0,1] = 0
$arr[0,2] = 10999
$arr[0,3] = 0
$arr[0,4] = 10999
$arr[0,5] = ""
$arr[0,6] = -11589
$ws.Range("H22:N22").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10999
$arr[0,1] = 0
$arr[0,2] = 10999
$arr[0,3] = 0
$arr[0,4] = 10999
$arr[0,5] = ""
$arr[0,6] = -11213
$ws.Range("H27:N27").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2549.7144
$arr[0,1] = 3499.5
$arr[0,2] = 2169.8
$arr[0,3] = 3499.5
$arr[0,4] = 2169.8
$arr[0,5] = -3311.5
$arr[0,6] = -2545.8
$ws.Range("H46:N46").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4473
$arr[0,1] = 4000
$arr[0,2] = 4630.6665
$arr[0,3] = 4000
$arr[0,4] = 4630.6665
$arr[0,5] = -3251
$arr[0,6] = -6128.6665
$ws.Range("H68:N68").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4473
$arr[0,1] = 4000
$arr[0,2] = 4630.6665
$arr[0,3] = 20000
$arr[0,4] = 23153.3325
$arr[0,5] = -16256
$arr[0,6] = -30641.3325
$ws.Range("H71:N71").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1415.8
$arr[0,1] = 464.3846
$arr[0,2] = 7600
$arr[0,3] = 1393.1538
$arr[0,4] = 22800
$arr[0,5] = 776.8462
$arr[0,6] = -27140
$ws.Range("H113:N113").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 73689
$arr[0,1] = 50000
$arr[0,2] = 97378
$arr[0,3] = 50000
$arr[0,4] = 97378
$arr[0,5] = -48433
$arr[0,6] = -100512
$ws.Range("H115:N115").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4230.0713
$arr[0,1] = 3523.3
$arr[0,2] = 5997
$arr[0,3] = 10569.9
$arr[0,4] = 17991
$arr[0,5] = -8119.900000000001
$arr[0,6] = -22891
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4019.889
$arr[0,1] = 4129.4614
$arr[0,2] = 3735
$arr[0,3] = 12388.3842
$arr[0,4] = 11205
$arr[0,5] = -9918.3842
$arr[0,6] = -16145
$ws.Range("H126:N126").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 90999.5
$arr[0,1] = 0
$arr[0,2] = 90999.5
$arr[0,3] = 0
$arr[0,4] = 90999.5
$arr[0,5] = ""
$arr[0,6] = -101119.5
$ws.Range("H133:N133").Value = $arr
